$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-06 Wednesday" "2024-11-07 Thursday"
Replace-Text "8+52=60" "92-82=10"
Replace-Text "83-20=63" "86-76=10"
Replace-Text "65-21=44" "0+16=16"
Replace-Text "1+85=86" "69-35=34"
Replace-Text "98-38=60" "83-60=23"
Replace-Text "97-1=96" "90-47=43"
Replace-Text "99-16=83" "44+22=66"
Replace-Text "19+58=77" "53-9=44"
Replace-Text "86-16=70" "99-12=87"
Replace-Text "57-2=55" "8+33=41"
Replace-Text "54+17=71" "87+6=93"
Replace-Text "86-57=29" "38+43=81"
Replace-Text "97-91=6" "96-14=82"
Replace-Text "17+49=66" "92-33=59"
Replace-Text "30-29=1" "19+32=51"
Replace-Text "16+19=35" "36-22=14"
Replace-Text "14+6=20" "45+50=95"
Replace-Text "66+27=93" "97-35=62"
Replace-Text "3+31=34" "27+45=72"
Replace-Text "17+42=59" "83-44=39"
Replace-Text "27-7=20" "94-90=4"
Replace-Text "39+42=81" "46+50=96"
Replace-Text "31+66=97" "6+17=23"
Replace-Text "29-24=5" "39+48=87"
Replace-Text "18-14=4" "88-9=79"
Replace-Text "0+2=2" "91+1=92"
Replace-Text "72-50=22" "57+26=83"
Replace-Text "59-37=22" "17+55=72"
Replace-Text "91-45=46" "15+49=64"
Replace-Text "45+43=88" "26+43=69"
Replace-Text "27+55=82" "79-59=20"
Replace-Text "56-39=17" "15+59=74"
Replace-Text "20+39=59" "75-15=60"
Replace-Text "23+62=85" "10+70=80"
Replace-Text "58+17=75" "87-78=9"
Replace-Text "25+45=70" "28+67=95"
Replace-Text "60+11=71" "67+9=76"
Replace-Text "29+60=89" "79-47=32"
Replace-Text "10+8=18" "80-50=30"
Replace-Text "78+7=85" "15+36=51"
Replace-Text "20+68=88" "3+63=66"
Replace-Text "46+26=72" "98+1=99"
Replace-Text "84-59=25" "58-12=46"
Replace-Text "28+0=28" "10+85=95"
Replace-Text "79-58=21" "44-9=35"
Replace-Text "75-0=75" "67-37=30"
Replace-Text "36+7=43" "36+35=71"
Replace-Text "42+32=74" "92-22=70"
Replace-Text "76-63=13" "31+52=83"
Replace-Text "77-29=48" "43+21=64"
Replace-Text "34-27=7" "28+67=95"
Replace-Text "46+37=83" "67-17=50"
Replace-Text "22+69=91" "62-56=6"
Replace-Text "7+32=39" "23-7=16"
Replace-Text "89-9=80" "45+36=81"
Replace-Text "31+55=86" "46+4=50"
Replace-Text "30+32=62" "28+24=52"
Replace-Text "83-12=71" "90-20=70"
Replace-Text "4+43=47" "80+12=92"
Replace-Text "31-4=27" "0+9=9"
Replace-Text "41+16=57" "26-14=12"
Replace-Text "57-55=2" "8+49=57"
Replace-Text "21+38=59" "18+58=76"
Replace-Text "42+4=46" "31+50=81"
Replace-Text "46-36=10" "59-41=18"
Replace-Text "96-34=62" "58-57=1"
Replace-Text "54-9=45" "77+12=89"
Replace-Text "54+36=90" "11+39=50"
Replace-Text "5+41=46" "62-45=17"
Replace-Text "55+17=72" "46-8=38"
Replace-Text "22+32=54" "80+11=91"
Replace-Text "38-16=22" "52+24=76"
Replace-Text "7+38=45" "10+39=49"
Replace-Text "35+52=87" "95-72=23"
Replace-Text "50-33=17" "47+28=75"
Replace-Text "96-35=61" "44+47=91"
Replace-Text "46+31=77" "62-18=44"
Replace-Text "61-12=49" "72+2=74"
Replace-Text "43-19=24" "44-27=17"
Replace-Text "48-4=44" "13+72=85"
Replace-Text "33-26=7" "55-47=8"
Replace-Text "51-48=3" "88-61=27"
Replace-Text "80-72=8" "7+70=77"
Replace-Text "31+8=39" "83-35=48"
Replace-Text "98-14=84" "2+19=21"
Replace-Text "61-14=47" "22-8=14"
Replace-Text "19+22=41" "3+66=69"
Replace-Text "50+28=78" "37-37=0"
Replace-Text "18+8=26" "46+50=96"
Replace-Text "26+11=37" "2+79=81"
Replace-Text "36+23=59" "89-4=85"
Replace-Text "15+45=60" "76-56=20"
Replace-Text "88-88=0" "55-19=36"
Replace-Text "81-42=39" "4+69=73"
Replace-Text "39+27=66" "27+55=82"
Replace-Text "97-53=44" "28-14=14"
Replace-Text "32-13=19" "5+68=73"
Replace-Text "6+78=84" "85-3=82"
Replace-Text "36-17=19" "9+21=30"
Replace-Text "53-44=9" "9+68=77"
